$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Resolution changes from "Unresolved" to "Fixed"
$ws.Cells.Item(10, 5).Value = "Fixed"

# Row 11: Rank changes from 10 to 11
$ws.Cells.Item(11, 7).Value = 11

# Row 12: Rank changes from 11 to 12
$ws.Cells.Item(12, 7).Value = 12

# Row 13: add Sprint ("Sprint 5") in column F, and change Rank from 12 to 10
$ws.Cells.Item(13, 6).Value = "Sprint 5"
$ws.Cells.Item(13, 7).Value = 10

# New Row 17: new PBI about adding/deleting sprints
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Als PO möchte ich Sprints hinzufügen und löschen können."
$ws.Cells.Item(17, 3).Value = "Akzeptanzkriterien:`nDie Sprints werden nach aufsteigendem Start-Datum sortiert"
$ws.Cells.Item(17, 4).Value = 3
$ws.Cells.Item(17, 5).Value = "Unresolved"
$ws.Cells.Item(17, 7).Value = 16

# Copy description cell formatting (wrap text style) from an existing description cell
$ws.Cells.Item(13, 3).Copy() | Out-Null
$ws.Cells.Item(17, 3).PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(17).RowHeight = 60

# Update selection to match new focus cell
$ws.Activate()
$ws.Range("F13").Select() | Out-Null
